# Apply the dated-worksheet update:
#  - bump the date heading by one day
#  - replace each "A×B=" problem in the practice table with a new one,
#    addressed by (row, column) so that duplicate problems (e.g. the two
#    "332×7=" cells) are updated independently to their correct targets.

$d = $word.ActiveDocument

# 1) Update the date heading paragraph.
$d.Content.Find.Execute(
    "2024-09-21 Saturday", $true, $false, $false, $false, $false,
    $true, 1, $false, "2024-09-22 Sunday", 2) | Out-Null

# 2) Update the multiplication-problem table, cell by cell.
$table = $d.Tables.Item(1)

function Set-CellText($row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $range = $cell.Range
    # Trim the trailing cell-mark character so we don't clobber it.
    $range.MoveEnd(1, -1) | Out-Null
    $range.Text = $newText
}

Set-CellText 1 1 "959×3="
Set-CellText 1 2 "687×3="
Set-CellText 1 3 "267×9="
Set-CellText 1 4 "862×5="
Set-CellText 1 5 "163×9="

Set-CellText 5 1 "829×2="
Set-CellText 5 2 "754×9="
Set-CellText 5 3 "598×3="
Set-CellText 5 4 "422×4="
Set-CellText 5 5 "432×9="

Set-CellText 10 1 "310×4="
Set-CellText 10 2 "385×8="
Set-CellText 10 3 "142×2="
Set-CellText 10 4 "687×4="
Set-CellText 10 5 "111×9="

Set-CellText 15 1 "196×6="
Set-CellText 15 2 "611×3="
Set-CellText 15 3 "131×4="
Set-CellText 15 4 "597×7="
Set-CellText 15 5 "251×9="

Set-CellText 20 1 "991×6="
Set-CellText 20 2 "601×3="
Set-CellText 20 3 "313×6="
Set-CellText 20 4 "234×3="
Set-CellText 20 5 "868×9="

Write-Output "Done updating date and table cells."
